$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 15: the "_" placeholder values in the UPI row get corrected to " UPI"
$ws.Range("B15").Value = " UPI"
$ws.Range("C15").Value = " UPI"

# Row 16: new UPI / CREDIT / RUPAY scheme row, copied (format + all) from row 15
$ws.Range("A15:H15").Copy()
$ws.Range("A16:H16").PasteSpecial(-4122)

$ws.Range("A16").Value = " UPI"
$ws.Range("B16").Value = "CREDIT"
$ws.Range("C16").Value = "RUPAY"
$ws.Range("D16").Value = 1.5
$ws.Range("E16").Value = 1.1
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = "Active / Inactive"
$ws.Range("H16").Value = "Start date of scheme"

$ws.Range("D6").Select()
